$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target columns remain text so numeric-looking values
# (e.g. "1.002", "29.411.66") are not auto-converted to numbers/dates by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.411.66'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '1.901.73'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').Value = '325.69'
$ws.Range('E5').Value = '  -2.89%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').Value = '0.4805'
$ws.Range('E7').Value = '  +2.70%  '
$ws.Range('D8').Value = '0.4066'
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('D9').Value = '0.08075'
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('D10').Value = '1.003'
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('D11').Value = '23.26'
$ws.Range('E11').Value = '  +3.48%  '
$ws.Range('D12').Value = '1.929.03'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').Value = '5.949'
$ws.Range('E13').Value = '  -0.78%  '
$ws.Range('D14').Value = '7.070'
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('D15').Value = '89.88'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').Value = '0.06699'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('D18').Value = '0.00001030'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').Value = '17.65'
$ws.Range('E19').Value = '  -1.25%  '
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').Value = '29.415.79'
$ws.Range('E21').Value = '  -0.97%  '
$ws.Range('D22').Value = '5.533'
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('D23').Value = '11.73'
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('D25').Value = '2.151.20'
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('D26').Value = '154.90'
$ws.Range('E26').Value = '  -0.85%  '
$ws.Range('D27').Value = '19.77'
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('D28').Value = '6.086'
$ws.Range('D29').Value = '2.093'
$ws.Range('E29').Value = '  -2.55%  '
$ws.Range('D30').Value = '118.30'
$ws.Range('E30').Value = '  +0.60%  '
$ws.Range('D31').Value = '1.025'
$ws.Range('E31').Value = '  -4.07%  '
$ws.Range('D32').Value = '0.09519'
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('D33').Value = '1.392'
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '5.394'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '3.529'
$ws.Range('E35').Value = '  -1.12%  '
$ws.Range('D36').Value = '0.02254'
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('D37').Value = '0.06060'
$ws.Range('E37').Value = '  -1.31%  '
$ws.Range('D38').Value = '1.173'
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('D39').Value = '0.5878'
$ws.Range('E39').Value = '  -0.33%  '
$ws.Range('D40').Value = '7.875'
$ws.Range('E40').Value = '  -6.88%  '
$ws.Range('D41').Value = '0.1845'
$ws.Range('D42').Value = '10.22'
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('D43').Value = '1.284'
$ws.Range('E43').Value = '  +1.94%  '
$ws.Range('D44').Value = '2.408'
$ws.Range('E44').Value = '  +1.94%  '
$ws.Range('D45').Value = '0.07738'
$ws.Range('E45').Value = '  +3.01%  '
$ws.Range('D46').Value = '12.29'
$ws.Range('E46').Value = '  +0.60%  '
$ws.Range('D47').Value = '0.5520'
$ws.Range('E47').Value = '  -1.14%  '
$ws.Range('D48').Value = '1.924'
$ws.Range('E48').Value = '  -0.67%  '
$ws.Range('D49').Value = '113.26'
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('D50').Value = '72.16'
$ws.Range('E50').Value = '  +0.92%  '
$ws.Range('D51').Value = '0.2927'
$ws.Range('E51').Value = '  -2.99%  '
